# Foot Pedal main functions update — now triggers tension.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2-11 (the "Row 2" foot-pedal function table) -----------------
# Column C ("Row 0" / scene+track block): a few labels were renamed/shifted.
$ws.Range("C5").Value  = "Tension Start"
$ws.Range("C6").Value  = "Tension Cancel"
$ws.Range("C7").Value  = "Go To Tension"
$ws.Range("C8").Value  = "Disarm All Tracks"
$ws.Range("C9").Value  = "Disable Clips of Armed Tracks"
$ws.Range("C10").Value = "Disable Notes of Armed Tracks"

# Column E ("Row 2" / chord-player block) is replaced by the foot-pedal /
# glove mapping commands.
$ws.Range("E2").Value  = "Map Foot Pedal 1"
$ws.Range("E3").Value  = "Map Foot Pedal 2"
$ws.Range("E4").Value  = "Map Glove X"
$ws.Range("E5").Value  = "Map Glove Y"
$ws.Range("E6").Value  = "Map Glove Z"
$ws.Range("E7").Value  = "UnMap Foot Pedal 1"
$ws.Range("E8").Value  = "UnMap Foot Pedal 2"
$ws.Range("E9").Value  = "Unmap Glove X"
$ws.Range("E10").Value = "Unmap Glove Y"
$ws.Range("E11").Value = "Unmap Glove Z"

# Columns F, G, H ("Row 3"/"Row 4"/"Row 5" blocks) are no longer used for
# the foot-pedal assignment — they now just read "Nothing".
$ws.Range("F2").Value  = "Nothing"
$ws.Range("F3").Value  = "Nothing"
$ws.Range("G2").Value  = "Nothing"
$ws.Range("G3").Value  = "Nothing"
$ws.Range("G4").Value  = "Nothing"
$ws.Range("G5").Value  = "Nothing"
$ws.Range("G6").Value  = "Nothing"
$ws.Range("G7").Value  = "Nothing"
$ws.Range("G8").Value  = "Nothing"
$ws.Range("G9").Value  = "Nothing"
$ws.Range("G10").Value = "Nothing"
$ws.Range("H2").Value  = "Nothing"
$ws.Range("H3").Value  = "Nothing"
$ws.Range("H4").Value  = "Nothing"
$ws.Range("H5").Value  = "Nothing"
$ws.Range("F7").Value  = "Nothing"
$ws.Range("F8").Value  = "Nothing"

# --- Row 16-19 (the "Row 8" block): now shows the Fishman hold state ---
$ws.Range("E16").Value = "Fishman Hold Unknown"
$ws.Range("E17").Value = "Fishman Hold Unknown"
$ws.Range("E18").Value = "Fishman Hold Unknown"
$ws.Range("E19").Value = "Fishman Hold Unknown"

# The old chord-player / glove-serial / midi-learn labels in column E used
# a red font (style index 1); that highlight no longer applies to the new
# foot-pedal/glove-mapping text, so clear it back to automatic/black.
$ws.Range("E2:E11").Font.Color = 0

# Selection moves to C2 in the saved file.
$ws.Range("C2").Select()
